$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so purely-numeric-looking values
# (e.g. "306.87") are not silently converted to numbers, matching the
# original inlineStr cell type used throughout column D.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.653.58'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '2.418.43'
$ws.Range("E3").Value = '  +2.41%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '306.87'
$ws.Range("D6").Value = '97.58'
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").Value = '35.05'
$ws.Range("E10").Value = '  +2.75%  '
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("E12").Value = '  +2.81%  '
$ws.Range("D13").Value = '18.56'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("D15").Value = '2.784.30'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '2.387.24'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("E17").Value = '  +3.59%  '
$ws.Range("D18").Value = '43.633.34'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").Value = '6.42'
$ws.Range("E19").Value = '  +2.16%  '
$ws.Range("D20").Value = '12.13'
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").Value = '68.30'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '238.69'
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '24.99'
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("E28").Value = '  -0.62%  '
$ws.Range("E29").Value = '  +3.30%  '
$ws.Range("D30").Value = '32.41'
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("D31").Value = '0.120'
$ws.Range("E31").Value = '  +18.47%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.14'
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("B33").Value = 'Celestia'
$ws.Range("C33").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D33").Value = '18.40'
$ws.Range("E33").Value = '  +6.90%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '0.0753'
$ws.Range("E35").Value = '  +3.47%  '
$ws.Range("D36").Value = '1.90'
$ws.Range("E36").Value = '  +3.53%  '
$ws.Range("D37").Value = '130.51'
$ws.Range("E37").Value = '  +26.17%  '
$ws.Range("D38").Value = '2.91'
$ws.Range("E38").Value = '  +5.30%  '
$ws.Range("D39").Value = '4.40'
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("D41").Value = '0.109'
$ws.Range("E41").Value = '  +0.19%  '
$ws.Range("D42").Value = '21.17'
$ws.Range("E42").Value = '  -5.87%  '
$ws.Range("D43").Value = '1.948.30'
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").Value = '0.0284'
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("D46").Value = '2.83'
$ws.Range("E46").Value = '  +3.35%  '
$ws.Range("D47").Value = '9.32'
$ws.Range("E47").Value = '  -1.45%  '
$ws.Range("D48").Value = '2.637.48'
$ws.Range("E48").Value = '  +1.90%  '
$ws.Range("E49").Value = '  +3.86%  '
$ws.Range("D50").Value = '52.86'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").Value = '72.40'
$ws.Range("E51").Value = '  +0.07%  '
